# Update "想去人数" (want-to-go count) values in column F for the two
# sheets that carry the real data table: "展览" and "全部类型".
# Both sheets mirror each other, so the same row/value updates apply to both.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 1533
    "F4"  = 990
    "F6"  = 2448
    "F8"  = 1509
    "F9"  = 71
    "F10" = 178
    "F11" = 55
    "F12" = 443
    "F14" = 23
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
